$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Meeting" -> "Meeting  " (two trailing spaces) in D3
$ws.Range("D3").Value = "Meeting  "

# "2 uur" -> "1/2 uur" in C3
$ws.Range("C3").Value = "1/2 uur"

# Move the active selection from H11 to C4
$ws.Activate()
$ws.Range("C4").Select()
